$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This quarterly report drops the oldest quarter (column D) and appends a new
# quarter (column M). Shift all the per-quarter data (labels, publish dates,
# figures) one column to the left (E->D, F->E, ... M->L), then populate the
# now-vacant column M with the new quarter's figures.
# Range.Copy is used (instead of Value2=Value2) so that text which happens to
# look like a date (e.g. "1401-04-30") is carried over verbatim as text
# instead of being re-interpreted/coerced into a date serial number.
# ---------------------------------------------------------------------------

$ws.Range("E8:M9").Copy($ws.Range("D8:L9"))
$ws.Range("E11:M27").Copy($ws.Range("D11:L27"))

# New quarter label + publish date (column M).
$ws.Range("M8").Value2 = "فصل چهارم منتهی به 1401/12"

# "1402-02-28" looks like a date, so build it as a text formula result in a
# scratch cell and paste-special (values only) into M9 to keep it as text
# without disturbing M9's existing cell style.
$ws.Range("ZZ1").Formula = '="1402" & "-02-28"'
$ws.Range("ZZ1").Copy()
$ws.Range("M9").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# New quarter financial figures (column M).
$ws.Range("M11").Value2 = 5977172
$ws.Range("M12").Value2 = -3400439
$ws.Range("M13").Value2 = 2576733
$ws.Range("M14").Value2 = -126788
$ws.Range("M15").Value2 = 0
$ws.Range("M16").Value2 = -163817
$ws.Range("M17").Value2 = 2286128
$ws.Range("M18").Value2 = -542924
$ws.Range("M19").Value2 = 335738
$ws.Range("M20").Value2 = 2078942
$ws.Range("M21").Value2 = 325229
$ws.Range("M22").Value2 = 2404171
$ws.Range("M23").Value2 = 0
$ws.Range("M24").Value2 = 2404171
$ws.Range("M25").Value2 = 365
$ws.Range("M26").Value2 = 6580000
$ws.Range("M27").Value2 = 365
